# Update September (R/S) and August (P/Q) transaction logs on sheet "2024":
# a new September entry was logged, pushing prior September rows down by one,
# and a new August entry was logged, pushing prior August rows down by one,
# which also pushes the trailing "Broadband" label from A62 down to A63.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# --- September_Details / September_Date (columns R/S), rows 30-58 ---
$ws.Range("R30").Value = "axis bna"
$ws.Range("S30").Value = "2024-09-05 14:18:32"
$ws.Range("R31").Value = "axis bna"
$ws.Range("S31").Value = "2024-09-05 14:13:16"
$ws.Range("R32").Value = "axis bna"
$ws.Range("S32").Value = "2024-09-05 14:15:23"
$ws.Range("R33").Value = "balance your axis"
$ws.Range("S33").Value = "2024-09-05 09:20:57"
$ws.Range("R34").Value = "bal axis"
$ws.Range("S34").Value = "2024-09-05 09:06:25"
$ws.Range("R35").Value = "broker"
$ws.Range("S35").Value = "2024-09-04 21:20:47"
$ws.Range("R36").Value = "exclusive on axis"
$ws.Range("S36").Value = "2024-09-04 13:21:05"
$ws.Range("R37").Value = "your corporate axis"
$ws.Range("S37").Value = "2024-09-04 11:46:10"
$ws.Range("R38").Value = "balance your axis"
$ws.Range("S38").Value = "2024-09-04 08:14:16"
$ws.Range("R39").Value = "axis"
$ws.Range("S39").Value = "2024-09-04 07:02:13"
$ws.Range("R40").Value = "bal axisbank w axis"
$ws.Range("S40").Value = "2024-09-04 06:53:15"
$ws.Range("R41").Value = "logging iob internet"
$ws.Range("S41").Value = "2024-09-03 20:09:12"
$ws.Range("R42").Value = "password internet"
$ws.Range("S42").Value = "2024-09-03 20:05:31"
$ws.Range("R43").Value = "logging iob internet"
$ws.Range("S43").Value = "2024-09-03 20:05:09"
$ws.Range("R44").Value = "internet"
$ws.Range("S44").Value = "2024-09-03 19:58:18"
$ws.Range("R45").Value = "login internet invalid"
$ws.Range("S45").Value = "2024-09-03 19:54:49"
$ws.Range("R46").Value = "login internet invalid"
$ws.Range("S46").Value = "2024-09-03 19:56:17"
$ws.Range("R47").Value = "corporate internet share"
$ws.Range("S47").Value = "2024-09-03 19:22:58"
$ws.Range("R48").Value = "login sbi internet personal do not share anyone"
$ws.Range("S48").Value = "2024-09-03 19:17:10"
$ws.Range("R49").Value = "login internet personal share"
$ws.Range("S49").Value = "2024-09-03 19:13:40"
$ws.Range("R50").Value = "internet verify it"
$ws.Range("S50").Value = "2024-09-03 19:05:49"
$ws.Range("R51").Value = "balance your axis"
$ws.Range("S51").Value = "2024-09-03 13:14:06"
$ws.Range("R52").Value = "lounge"
$ws.Range("S52").Value = "2024-09-03 13:08:08"
$ws.Range("R53").Value = "balance your axis"
$ws.Range("S53").Value = "2024-09-03 11:21:30"
$ws.Range("R54").Value = "broker"
$ws.Range("S54").Value = "2024-09-01 22:35:38"
$ws.Range("R55").Value = "amazeloan"
$ws.Range("S55").Value = "2024-09-01 10:12:03"
$ws.Range("R56").Value = "amazeloan"
$ws.Range("S56").Value = "2024-09-01 09:42:38"
$ws.Range("R57").Value = "amazeloan"
$ws.Range("S57").Value = "2024-09-01 09:29:24"
$ws.Range("R58").Value = "amazeloan"
$ws.Range("S58").Value = "2024-09-01 09:27:06"

# --- August_Details / August_Date (columns P/Q), rows 58-62 ---
$ws.Range("P58").Value = ""
$ws.Range("Q58").Value = ""
$ws.Range("P59").Value = "hdfc"
$ws.Range("Q59").Value = "2024-08-30 12:15:48"
$ws.Range("P60").Value = "hdfc"
$ws.Range("Q60").Value = "2024-08-21 20:17:10"
$ws.Range("P61").Value = "hdfc"
$ws.Range("Q61").Value = "2024-08-21 20:16:45"
$ws.Range("P62").Value = "hdfc"
$ws.Range("Q62").Value = "2024-08-21 20:15:50"

# --- "Broadband" label moves from row 62 to the newly appended row 63 ---
$ws.Range("A62").Value = ""
$ws.Range("A63").Value = "Broadband"

